$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data block (row 1 header stays untouched, keeping its bold/border style
# and leaving styles.xml unmodified). This also drops the now-unused shared strings so that
# we can rebuild the shared-string table in the exact order the target file expects.
$ws.Range("A2:T13").Clear()

# Rewrite the text columns (A-D) column-by-column (top-to-bottom within each column, then
# moving to the next column) so new distinct strings get interned into the shared-string
# table in this exact order: ECs, FAPs, M2, sCs, Adam23, Itga5.

# Column A
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(17, 1).Value = "sCs"

# Column B
$ws.Cells.Item(2, 2).Value = "Adam23"
$ws.Cells.Item(3, 2).Value = "Adam23"
$ws.Cells.Item(4, 2).Value = "Adam23"
$ws.Cells.Item(5, 2).Value = "Adam23"
$ws.Cells.Item(6, 2).Value = "Adam23"
$ws.Cells.Item(7, 2).Value = "Adam23"
$ws.Cells.Item(8, 2).Value = "Adam23"
$ws.Cells.Item(9, 2).Value = "Adam23"
$ws.Cells.Item(10, 2).Value = "Adam23"
$ws.Cells.Item(11, 2).Value = "Adam23"
$ws.Cells.Item(12, 2).Value = "Adam23"
$ws.Cells.Item(13, 2).Value = "Adam23"
$ws.Cells.Item(14, 2).Value = "Adam23"
$ws.Cells.Item(15, 2).Value = "Adam23"
$ws.Cells.Item(16, 2).Value = "Adam23"
$ws.Cells.Item(17, 2).Value = "Adam23"

# Column C
$ws.Cells.Item(2, 3).Value = "Itga5"
$ws.Cells.Item(3, 3).Value = "Itga5"
$ws.Cells.Item(4, 3).Value = "Itga5"
$ws.Cells.Item(5, 3).Value = "Itga5"
$ws.Cells.Item(6, 3).Value = "Itga5"
$ws.Cells.Item(7, 3).Value = "Itga5"
$ws.Cells.Item(8, 3).Value = "Itga5"
$ws.Cells.Item(9, 3).Value = "Itga5"
$ws.Cells.Item(10, 3).Value = "Itga5"
$ws.Cells.Item(11, 3).Value = "Itga5"
$ws.Cells.Item(12, 3).Value = "Itga5"
$ws.Cells.Item(13, 3).Value = "Itga5"
$ws.Cells.Item(14, 3).Value = "Itga5"
$ws.Cells.Item(15, 3).Value = "Itga5"
$ws.Cells.Item(16, 3).Value = "Itga5"
$ws.Cells.Item(17, 3).Value = "Itga5"

# Column D
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(17, 4).Value = "sCs"

# Now fill in the numeric columns (E-T) for every row; order amongst these does not
# affect the shared-string table since they hold numbers, not text.

# Row 2 numeric values
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.269164
$ws.Cells.Item(2, 8).Value = 0.8074920000000001
$ws.Cells.Item(2, 9).Value = 0.01089095165781685
$ws.Cells.Item(2, 10).Value = 0.01089095165781686
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 30.99161333333333
$ws.Cells.Item(2, 14).Value = 92.97484
$ws.Cells.Item(2, 15).Value = 0.3599121977633812
$ws.Cells.Item(2, 16).Value = 0.3599121977633811
$ws.Cells.Item(2, 17).Value = 8.341826611253333
$ws.Cells.Item(2, 18).Value = 75.07643950128001
$ws.Cells.Item(2, 19).Value = 0.003919786346899604
$ws.Cells.Item(2, 20).Value = 0.003919786346899603

# Row 3 numeric values
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.269164
$ws.Cells.Item(3, 8).Value = 0.8074920000000001
$ws.Cells.Item(3, 9).Value = 0.01089095165781685
$ws.Cells.Item(3, 10).Value = 0.01089095165781686
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.913269
$ws.Cells.Item(3, 14).Value = 89.739807
$ws.Cells.Item(3, 15).Value = 0.3473891556493311
$ws.Cells.Item(3, 16).Value = 0.3473891556493311
$ws.Cells.Item(3, 17).Value = 8.051575137116
$ws.Cells.Item(3, 18).Value = 72.464176234044
$ws.Cells.Item(3, 19).Value = 0.00378339850062668
$ws.Cells.Item(3, 20).Value = 0.00378339850062668

# Row 4 numeric values
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.269164
$ws.Cells.Item(4, 8).Value = 0.8074920000000001
$ws.Cells.Item(4, 9).Value = 0.01089095165781685
$ws.Cells.Item(4, 10).Value = 0.01089095165781686
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 19.150218
$ws.Cells.Item(4, 14).Value = 57.450654
$ws.Cells.Item(4, 15).Value = 0.2223955550134164
$ws.Cells.Item(4, 16).Value = 0.2223955550134163
$ws.Cells.Item(4, 17).Value = 5.154549277752
$ws.Cells.Item(4, 18).Value = 46.39094349976801
$ws.Cells.Item(4, 19).Value = 0.002422099238564467
$ws.Cells.Item(4, 20).Value = 0.002422099238564466

# Row 5 numeric values
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.269164
$ws.Cells.Item(5, 8).Value = 0.8074920000000001
$ws.Cells.Item(5, 9).Value = 0.01089095165781685
$ws.Cells.Item(5, 10).Value = 0.01089095165781686
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.053716000000001
$ws.Cells.Item(5, 14).Value = 18.161148
$ws.Cells.Item(5, 15).Value = 0.07030309157387134
$ws.Cells.Item(5, 16).Value = 0.07030309157387132
$ws.Cells.Item(5, 17).Value = 1.629442413424
$ws.Cells.Item(5, 18).Value = 14.664981720816
$ws.Cells.Item(5, 19).Value = 0.0007656675717261041
$ws.Cells.Item(5, 20).Value = 0.0007656675717261041

# Row 6 numeric values
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 20.07911633333334
$ws.Cells.Item(6, 8).Value = 60.237349
$ws.Cells.Item(6, 9).Value = 0.8124440315867432
$ws.Cells.Item(6, 10).Value = 0.8124440315867433
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 30.99161333333333
$ws.Cells.Item(6, 14).Value = 92.97484
$ws.Cells.Item(6, 15).Value = 0.3599121977633812
$ws.Cells.Item(6, 16).Value = 0.3599121977633811
$ws.Cells.Item(6, 17).Value = 622.2842094776845
$ws.Cells.Item(6, 18).Value = 5600.55788529916
$ws.Cells.Item(6, 19).Value = 0.2924085169681266
$ws.Cells.Item(6, 20).Value = 0.2924085169681266

# Row 7 numeric values
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 20.07911633333334
$ws.Cells.Item(7, 8).Value = 60.237349
$ws.Cells.Item(7, 9).Value = 0.8124440315867432
$ws.Cells.Item(7, 10).Value = 0.8124440315867433
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 29.913269
$ws.Cells.Item(7, 14).Value = 89.739807
$ws.Cells.Item(7, 15).Value = 0.3473891556493311
$ws.Cells.Item(7, 16).Value = 0.3473891556493311
$ws.Cells.Item(7, 17).Value = 600.6320081612937
$ws.Cells.Item(7, 18).Value = 5405.688073451644
$ws.Cells.Item(7, 19).Value = 0.2822342461452572
$ws.Cells.Item(7, 20).Value = 0.2822342461452572

# Row 8 numeric values
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 20.07911633333334
$ws.Cells.Item(8, 8).Value = 60.237349
$ws.Cells.Item(8, 9).Value = 0.8124440315867432
$ws.Cells.Item(8, 10).Value = 0.8124440315867433
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 19.150218
$ws.Cells.Item(8, 14).Value = 57.450654
$ws.Cells.Item(8, 15).Value = 0.2223955550134164
$ws.Cells.Item(8, 16).Value = 0.2223955550134163
$ws.Cells.Item(8, 17).Value = 384.519455030694
$ws.Cells.Item(8, 18).Value = 3460.675095276246
$ws.Cells.Item(8, 19).Value = 0.1806839413220714
$ws.Cells.Item(8, 20).Value = 0.1806839413220713

# Row 9 numeric values
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 20.07911633333334
$ws.Cells.Item(9, 8).Value = 60.237349
$ws.Cells.Item(9, 9).Value = 0.8124440315867432
$ws.Cells.Item(9, 10).Value = 0.8124440315867433
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 6.053716000000001
$ws.Cells.Item(9, 14).Value = 18.161148
$ws.Cells.Item(9, 15).Value = 0.07030309157387134
$ws.Cells.Item(9, 16).Value = 0.07030309157387132
$ws.Cells.Item(9, 17).Value = 121.5532678129614
$ws.Cells.Item(9, 18).Value = 1093.979410316652
$ws.Cells.Item(9, 19).Value = 0.05711732715128802
$ws.Cells.Item(9, 20).Value = 0.05711732715128802

# Row 10 numeric values
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.174593
$ws.Cells.Item(10, 8).Value = 0.523779
$ws.Cells.Item(10, 9).Value = 0.007064406543197522
$ws.Cells.Item(10, 10).Value = 0.007064406543197523
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 30.99161333333333
$ws.Cells.Item(10, 14).Value = 92.97484
$ws.Cells.Item(10, 15).Value = 0.3599121977633812
$ws.Cells.Item(10, 16).Value = 0.3599121977633811
$ws.Cells.Item(10, 17).Value = 5.410918746706667
$ws.Cells.Item(10, 18).Value = 48.69826872036
$ws.Cells.Item(10, 19).Value = 0.002542566084856231
$ws.Cells.Item(10, 20).Value = 0.002542566084856231

# Row 11 numeric values
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.174593
$ws.Cells.Item(11, 8).Value = 0.523779
$ws.Cells.Item(11, 9).Value = 0.007064406543197522
$ws.Cells.Item(11, 10).Value = 0.007064406543197523
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 29.913269
$ws.Cells.Item(11, 14).Value = 89.739807
$ws.Cells.Item(11, 15).Value = 0.3473891556493311
$ws.Cells.Item(11, 16).Value = 0.3473891556493311
$ws.Cells.Item(11, 17).Value = 5.222647374517
$ws.Cells.Item(11, 18).Value = 47.003826370653
$ws.Cells.Item(11, 19).Value = 0.002454098224204997
$ws.Cells.Item(11, 20).Value = 0.002454098224204997

# Row 12 numeric values
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.174593
$ws.Cells.Item(12, 8).Value = 0.523779
$ws.Cells.Item(12, 9).Value = 0.007064406543197522
$ws.Cells.Item(12, 10).Value = 0.007064406543197523
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 19.150218
$ws.Cells.Item(12, 14).Value = 57.450654
$ws.Cells.Item(12, 15).Value = 0.2223955550134164
$ws.Cells.Item(12, 16).Value = 0.2223955550134163
$ws.Cells.Item(12, 17).Value = 3.343494011274
$ws.Cells.Item(12, 18).Value = 30.091446101466
$ws.Cells.Item(12, 19).Value = 0.001571092614014823
$ws.Cells.Item(12, 20).Value = 0.001571092614014823

# Row 13 numeric values
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.174593
$ws.Cells.Item(13, 8).Value = 0.523779
$ws.Cells.Item(13, 9).Value = 0.007064406543197522
$ws.Cells.Item(13, 10).Value = 0.007064406543197523
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 6.053716000000001
$ws.Cells.Item(13, 14).Value = 18.161148
$ws.Cells.Item(13, 15).Value = 0.07030309157387134
$ws.Cells.Item(13, 16).Value = 0.07030309157387132
$ws.Cells.Item(13, 17).Value = 1.056936437588
$ws.Cells.Item(13, 18).Value = 9.512427938292001
$ws.Cells.Item(13, 19).Value = 0.0004966496201214713
$ws.Cells.Item(13, 20).Value = 0.0004966496201214713

# Row 14 numeric values
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 4.191587666666666
$ws.Cells.Item(14, 8).Value = 12.574763
$ws.Cells.Item(14, 9).Value = 0.1696006102122423
$ws.Cells.Item(14, 10).Value = 0.1696006102122424
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 30.99161333333333
$ws.Cells.Item(14, 14).Value = 92.97484
$ws.Cells.Item(14, 15).Value = 0.3599121977633812
$ws.Cells.Item(14, 16).Value = 0.3599121977633811
$ws.Cells.Item(14, 17).Value = 129.9040642181022
$ws.Cells.Item(14, 18).Value = 1169.13657796292
$ws.Cells.Item(14, 19).Value = 0.0610413283634987
$ws.Cells.Item(14, 20).Value = 0.0610413283634987

# Row 15 numeric values
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 4.191587666666666
$ws.Cells.Item(15, 8).Value = 12.574763
$ws.Cells.Item(15, 9).Value = 0.1696006102122423
$ws.Cells.Item(15, 10).Value = 0.1696006102122424
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 29.913269
$ws.Cells.Item(15, 14).Value = 89.739807
$ws.Cells.Item(15, 15).Value = 0.3473891556493311
$ws.Cells.Item(15, 16).Value = 0.3473891556493311
$ws.Cells.Item(15, 17).Value = 125.3840894100823
$ws.Cells.Item(15, 18).Value = 1128.456804690741
$ws.Cells.Item(15, 19).Value = 0.0589174127792422
$ws.Cells.Item(15, 20).Value = 0.05891741277924221

# Row 16 numeric values
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 4.191587666666666
$ws.Cells.Item(16, 8).Value = 12.574763
$ws.Cells.Item(16, 9).Value = 0.1696006102122423
$ws.Cells.Item(16, 10).Value = 0.1696006102122424
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 19.150218
$ws.Cells.Item(16, 14).Value = 57.450654
$ws.Cells.Item(16, 15).Value = 0.2223955550134164
$ws.Cells.Item(16, 16).Value = 0.2223955550134163
$ws.Cells.Item(16, 17).Value = 80.26981758277799
$ws.Cells.Item(16, 18).Value = 722.428358245002
$ws.Cells.Item(16, 19).Value = 0.03771842183876573
$ws.Cells.Item(16, 20).Value = 0.03771842183876574

# Row 17 numeric values
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 4.191587666666666
$ws.Cells.Item(17, 8).Value = 12.574763
$ws.Cells.Item(17, 9).Value = 0.1696006102122423
$ws.Cells.Item(17, 10).Value = 0.1696006102122424
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 6.053716000000001
$ws.Cells.Item(17, 14).Value = 18.161148
$ws.Cells.Item(17, 15).Value = 0.07030309157387134
$ws.Cells.Item(17, 16).Value = 0.07030309157387132
$ws.Cells.Item(17, 17).Value = 25.37468132310266
$ws.Cells.Item(17, 18).Value = 228.372131907924
$ws.Cells.Item(17, 19).Value = 0.01192344723073573
$ws.Cells.Item(17, 20).Value = 0.01192344723073573
